$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update token (C2)
$ws.Range("C2").Value = "be4ddf45ba067bd35faa80adfab8b45c"

# Update voto (F2)
$ws.Range("F2").Value = 4

# Update token (C3)
$ws.Range("C3").Value = "548830b810fdc6d8698e6600f557a067"

# Trim trailing whitespace from A33 correo
$ws.Range("A33").Value = "dcamerosv@miumg.edu.gt"

# Fill in missing token for C33
$ws.Range("C33").Value = "513de2189efa9ed6a7c1f618e582bb5a"
